# Profile new script implementation
# Adds a new test-case row (row 72: "Profile71") below the existing last
# row (row 71: "Profile70") on the "Test Cases" sheet, mirroring its
# formatting, and moves the sheet's viewport/selection down to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (A71:E71) into the new row (A72:E72)
# so the new row picks up the same cell styles (borders/fonts/wrap, etc.)
# as every other data row instead of Excel's default formatting.
$ws.Range("A71:E71").Copy($ws.Range("A72:E72"))

# Overwrite the copied values with the new test case's data.
$ws.Range("A72").Value = "Profile71"
$ws.Range("B72").Value = "OPQA-TBD"
$ws.Range("C72").Value = "Verify that First time logged user Profile tabs should display with default messags for each tab and count should be '0'"
# D72 (Runmode) and E72 (Results) stay as copied from row 71 ("Y" / blank).

# Scroll the view down and move the selection to just past the new last
# row, matching where a user would land after adding this row.
$win = $excel.ActiveWindow()
$win.ScrollRow() = 58
$win.ScrollColumn() = 1
$ws.Range("A73").Select()
